$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.366.22"
$ws.Range("E2").Value = "  +1.32%  "
$ws.Range("D3").Value = "1.886.98"
$ws.Range("E3").Value = "  +1.37%  "
$ws.Range("D4").Value = "'1.019"
$ws.Range("E4").Value = "  +1.47%  "
$ws.Range("D5").Value = "'316.78"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("E6").Value = "  +1.33%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'0.3935"
$ws.Range("D9").Value = "'0.08352"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("D10").Value = "'1.124"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("D11").Value = "'42.02"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'6.253"
$ws.Range("E12").Value = "  +0.85%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.897.45"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("D14").Value = "'20.50"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Value = "'7.287"
$ws.Range("E15").Value = "  +0.43%  "
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").Value = "'0.00001106"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "'90.91"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "'0.06711"
$ws.Range("E19").Value = "  +0.89%  "
$ws.Range("D20").Value = "'17.86"
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'1.017"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("D22").Value = "'6.052"
$ws.Range("D23").Value = "28.410.38"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").Value = "'11.18"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Value = "'2.306"
$ws.Range("E25").Value = "  +2.30%  "
$ws.Range("D26").Value = "2.115.60"
$ws.Range("E26").Value = "  +1.94%  "
$ws.Range("D27").Value = "'161.32"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "'2.455"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.72"
$ws.Range("E29").Value = "  +1.13%  "
$ws.Range("D30").Value = "'126.92"
$ws.Range("D31").Value = "'0.1060"
$ws.Range("E31").Value = "  -0.46%  "
$ws.Range("D32").Value = "'1.042"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").Value = "'5.913"
$ws.Range("E33").Value = "  -1.08%  "
$ws.Range("D34").Value = "'3.635"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("D35").Value = "'9.487"
$ws.Range("E35").Value = "  +1.42%  "
$ws.Range("D36").Value = "'0.02458"
$ws.Range("E36").Value = "  +1.70%  "
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("E38").Value = "  +2.58%  "
$ws.Range("D39").Value = "'0.6510"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "'1.254"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").Value = "'1.186"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").Value = "'5.017"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").Value = "'11.21"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").Value = "'0.6126"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("D45").Value = "'13.19"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("D46").Value = "'3.710"
$ws.Range("E46").Value = "  +1.39%  "
$ws.Range("D47").Value = "'1.288"
$ws.Range("E47").Value = "  +0.67%  "
$ws.Range("D48").Value = "'1.238"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("D49").Value = "'2.008"
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Value = "'121.47"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("D51").Value = "'0.06913"
$ws.Range("E51").Value = "  +1.08%  "
